# Remove the row for "ARMIJOS SALINAS LUIS CLAUDIO" (an all-zero row) from
# both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. All subsequent
# rows shift up by one, and the trailing summary row's "X de 59" counters
# become "X de 58" to reflect the reduced client count.

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Delete the data row (row 8) on each affected sheet; this shifts rows
# 9..61 up to 8..60, matching the diff exactly since the deleted row's
# values were all zero.
$wsGrupo.Rows.Item(8).Delete()
$wsMensual.Rows.Item(8).Delete()

# Fix up the "X de 59" -> "X de 58" summary counters on the last row of
# "VENTAS POR GRUPO" (now row 60, previously row 61), columns C..R.
$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R")
foreach ($col in $cols) {
    $cell = $wsGrupo.Range($col + "60")
    $cell.Value = ($cell.Value2 -replace "de 59", "de 58")
}
